$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lũy kế tháng LONG XUYÊN")

# Update the "last edited" timestamp text shared across the D column (rows 3,4,5,7,13)
$oldStamp = "2024-07-28T16:31:00.000Z"
$newStamp = "2024-07-31T18:24:00.000Z"
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldStamp) {
        $cell.Value = $newStamp
    }
}

# Update the T8 (Tháng 7) row numeric figures - row 5
$ws.Range("W5").Value = 43384000
$ws.Range("AA5").Value = 12566000
$ws.Range("AE5").Value = 55950000
$ws.Range("AH5").Value = 52650000
$ws.Range("AK5").Value = 15
$ws.Range("AQ5").Value = 58150000
